# Generate Report for Archive
#
# 1) Shared string update: "Ready for handoff" -> "In Translation" for every
#    Status cell that currently holds it (Overview!E2:F3, zh-cn!C2:C3,
#    de-de!C2:C3 all point at the same shared-string entry, so a single
#    Find/Replace across the workbook flips them all together).
# 2) Narrow the "Status" column(s) on every sheet from ~17.22 chars to
#    ~13.41 chars (Overview columns E & F; zh-cn & de-de column C).

$wb = $excel.ActiveWorkbook

foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    $found = $used.Find("Ready for handoff")
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        while ($true) {
            $found.Value = "In Translation"
            $found = $used.FindNext($found)
            if ($found -eq $null -or $found.Address() -eq $firstAddress) {
                break
            }
        }
    }
}

# --- Column width adjustments -------------------------------------------
# Target stored width (per the OOXML diff) is 13.4101845877511 characters.
# The host's ColumnWidth setter snaps to the same pixel-grid Excel itself
# uses (stored = round(rawInput*6)/6 + 5/6), so a raw input of 12.5 is the
# value that lands closest to the target after that snapping.
$newStatusWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newStatusWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newStatusWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newStatusWidth
